# Insert a new data row at row 61 (shifts existing rows 61-160 down to 62-161)
# and populate it with the new record, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(61).Insert()

$ws.Range("A61").Value = 11
$ws.Range("B61").Value = "Vega Monumental Concepción"
$ws.Range("C61").Value = "Bíobío"
$ws.Range("D61").Value = 44838
$ws.Range("E61").Value = 8
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100102
$ws.Range("H61").Value = "Cítricos"
$ws.Range("I61").Value = 100102004
$ws.Range("J61").Value = "Mandarina"
$ws.Range("K61").Value = "Murcott"
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 270
$ws.Range("N61").Value = 5000
$ws.Range("O61").Value = 5500
$ws.Range("P61").Value = 5278
$ws.Range("Q61").Value = "`$/bandeja 18 kilos"
$ws.Range("R61").Value = "Región de O'Higgins"
$ws.Range("S61").Value = 293
$ws.Range("T61").Value = 18
